$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Crystal Arcade": append rows 30-33
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Crystal Arcade")

# Template rows already present on the sheet to copy cell formatting from:
#   row 29 -> G style = "Equipo 1" (bold / blue fill)
#   row 28 -> G style = "Equipo 2" (bold / red fill)
$equipo1Template1 = "A29:N29"
$equipo2Template1 = "A28:N28"

$rows1 = @(
    @("JANET","FANG","TARA","LILY","DRACO","GUS","Equipo 1","TRB|Zeus 解開","TRB|Lxffy","TRB|R B M","NHG|Xemp","KCP|Fade","KCP|Tyrant","20250723T215318.000Z"),
    @("JANET","FANG","TARA","LILY","DRACO","GUS","Equipo 2","TRB|Zeus 解開","TRB|Lxffy","TRB|R B M","NHG|Xemp","KCP|Fade","KCP|Tyrant","20250723T215136.000Z"),
    @("JANET","FANG","TARA","LILY","DRACO","GUS","Equipo 1","TRB|Zeus 解開","TRB|Lxffy","TRB|R B M","NHG|Xemp","KCP|Fade","KCP|Tyrant","20250723T214926.000Z"),
    @("GRAY","FANG","EMZ","ALLI","SANDY","GUS","Equipo 2","TRB|Zeus 解開","TRB|Lxffy","TRB|R B M","NHG|Xemp","KCP|Fade","KCP|Tyrant","20250723T214405.000Z")
)

$startRow1 = 30
for ($i = 0; $i -lt $rows1.Length; $i++) {
    $r = $startRow1 + $i
    $data = $rows1[$i]
    $destRange = "A" + $r + ":N" + $r

    if ($data[6] -eq "Equipo 1") {
        $ws1.Range($equipo1Template1).Copy()
    } else {
        $ws1.Range($equipo2Template1).Copy()
    }
    $ws1.Range($destRange).PasteSpecial(-4122)

    for ($c = 1; $c -le $data.Length; $c++) {
        $ws1.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}

# ---------------------------------------------------------------------
# Sheet "Hot Potato": append rows 46-47
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Hot Potato")

# Template row for "Equipo 1" styling on this sheet
$equipo1Template2 = "A44:N44"

$rows2 = @(
    @("LUMI","DRACO","ALLI","BERRY","MICO","EMZ","Equipo 1","TRB|Zeus 解開","TRB|R B M","TRB|Lxffy","NHG|Xemp","KCP|Fade","KCP|Tyrant","20250723T220010.000Z"),
    @("LUMI","DRACO","ALLI","BERRY","MICO","EMZ","Equipo 1","TRB|Zeus 解開","TRB|R B M","TRB|Lxffy","NHG|Xemp","KCP|Fade","KCP|Tyrant","20250723T215842.000Z")
)

$startRow2 = 46
for ($i = 0; $i -lt $rows2.Length; $i++) {
    $r = $startRow2 + $i
    $data = $rows2[$i]
    $destRange = "A" + $r + ":N" + $r

    $ws2.Range($equipo1Template2).Copy()
    $ws2.Range($destRange).PasteSpecial(-4122)

    for ($c = 1; $c -le $data.Length; $c++) {
        $ws2.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}
